$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Banks(246)"
$ws.Range("B2").Value = 0.535864800889695
$ws.Range("A3").Value = "Airlines(14)"
$ws.Range("B3").Value = 0.4876383660172545
$ws.Range("A4").Value = "Auto Components(21)"
$ws.Range("B4").Value = 0.4735946512814284
$ws.Range("A5").Value = "Leisure Products(11)"
$ws.Range("B5").Value = 0.4685931722379817
$ws.Range("A6").Value = "Energy Equipment & Services(32)"
$ws.Range("B6").Value = 0.4588777949869692
$ws.Range("A7").Value = "Multi-Utilities(18)"
$ws.Range("B7").Value = 0.4491987284344434
$ws.Range("A8").Value = "Trading Companies & Distributors(25)"
$ws.Range("B8").Value = 0.4447503474417077
$ws.Range("A9").Value = "Machinery(85)"
$ws.Range("B9").Value = 0.420757989317279
$ws.Range("A10").Value = "Building Products(23)"
$ws.Range("B10").Value = 0.4164659188192756
$ws.Range("A11").Value = "Marine(15)"
$ws.Range("B11").Value = 0.4040424720855461
$ws.Range("A12").Value = "Road & Rail(22)"
$ws.Range("B12").Value = 0.372095315183784
$ws.Range("A13").Value = "Gas Utilities(12)"
$ws.Range("B13").Value = 0.3537754031430047
$ws.Range("A14").Value = "Semiconductors & Semiconductor Equipment(68)"
$ws.Range("B14").Value = 0.3526398566151605
$ws.Range("A15").Value = "Containers & Packaging(12)"
$ws.Range("B15").Value = 0.342474928558079
$ws.Range("A16").Value = "Chemicals(51)"
$ws.Range("B16").Value = 0.3364590752145379
$ws.Range("A17").Value = "Electric Utilities(28)"
$ws.Range("B17").Value = 0.335478499178457
$ws.Range("A18").Value = "Textiles, Apparel & Luxury Goods(29)"
$ws.Range("B18").Value = 0.3092906967350258
$ws.Range("A19").Value = "Oil, Gas & Consumable Fuels(122)"
$ws.Range("B19").Value = 0.3006234740052968
$ws.Range("A20").Value = "Thrifts & Mortgage Finance(47)"
$ws.Range("B20").Value = 0.2873668454502063
$ws.Range("A21").Value = "Construction & Engineering(20)"
$ws.Range("B21").Value = 0.2735771613150294
$ws.Range("A22").Value = "Life Sciences Tools & Services(19)"
$ws.Range("B22").Value = 0.266813741250306
$ws.Range("A23").Value = "Metals & Mining(89)"
$ws.Range("B23").Value = 0.2494195680278082
$ws.Range("A24").Value = "Insurance(75)"
$ws.Range("B24").Value = 0.2493363790682562
$ws.Range("A25").Value = "Electrical Equipment(28)"
$ws.Range("B25").Value = 0.2311967513597625
$ws.Range("A26").Value = "Specialty Retail(58)"
$ws.Range("B26").Value = 0.2298036480712669
$ws.Range("A27").Value = "Capital Markets(75)"
$ws.Range("B27").Value = 0.2231735188836178
$ws.Range("A28").Value = "Beverages(21)"
$ws.Range("B28").Value = 0.1984899453511643
$ws.Range("A29").Value = "Real Estate Management & Development(22)"
$ws.Range("B29").Value = 0.1832707141719118
$ws.Range("A30").Value = "Commercial Services & Supplies(52)"
$ws.Range("B30").Value = 0.1734925025275525
$ws.Range("A31").Value = "Hotels, Restaurants & Leisure(50)"
$ws.Range("B31").Value = 0.1622985598874865
$ws.Range("A32").Value = "Software(66)"
$ws.Range("B32").Value = 0.153694949467297
$ws.Range("A33").Value = "Health Care Providers & Services(46)"
$ws.Range("B33").Value = 0.1510486245465521
$ws.Range("A34").Value = "IT Services(52)"
$ws.Range("B34").Value = 0.1484281979832643
$ws.Range("A35").Value = "Household Durables(39)"
$ws.Range("B35").Value = 0.1409376841253368
$ws.Range("A36").Value = "Professional Services(35)"
$ws.Range("B36").Value = 0.1366943865403495
$ws.Range("A37").Value = "Food Products(44)"
$ws.Range("B37").Value = 0.1338765716751593
$ws.Range("A38").Value = "Pharmaceuticals(48)"
$ws.Range("B38").Value = 0.1311558671429781
$ws.Range("A39").Value = "Media(42)"
$ws.Range("B39").Value = 0.1142760206327241
$ws.Range("A40").Value = "Health Care Equipment & Supplies(83)"
$ws.Range("B40").Value = 0.1037151347855811
$ws.Range("A41").Value = "Communications Equipment(45)"
$ws.Range("B41").Value = 0.09691230739560371
$ws.Range("A42").Value = "Biotechnology(126)"
$ws.Range("B42").Value = 0.07326983662085039
